$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 68 appended to the historical price table.
# Force column A's new cell to be stored as plain text (matching the existing
# "Date" column cells) instead of letting Excel auto-convert the
# yyyy-mm-dd-looking string into a date serial number / date-formatted cell.
$ws.Cells.Item(68, 1).NumberFormat = "@"
$ws.Cells.Item(68, 1).Value = "2025-10-22"
$ws.Cells.Item(68, 1).ClearFormats()

$ws.Cells.Item(68, 2).Value = 54.18000030517578
$ws.Cells.Item(68, 3).Value = 401.8500061035156
$ws.Cells.Item(68, 4).Value = 338.1000061035156
